$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 15-18 down by one, with the old row 18
# wrapping around to become the new row 15:
#   new row 15 <- old row 18
#   new row 16 <- old row 15
#   new row 17 <- old row 16
#   new row 18 <- old row 17
#
# Capture every relevant cell value from the four rows first (values only
# change within columns A,B,E,F,G,H,P,Q,R,S,Z,AB,AC,AF,AW,AX - other columns
# are identical across the rotation so they don't need touching).

$cols = @("A","B","E","F","G","H","P","Q","R","S","Z","AB","AC","AF","AW","AX")

$data = @{}
foreach ($r in 15..18) {
    $data[$r] = @{}
    foreach ($c in $cols) {
        $data[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# new row -> source old row
$srcRow = @{ 15 = 18; 16 = 15; 17 = 16; 18 = 17 }

foreach ($r in 15..18) {
    $s = $srcRow[$r]

    $ws.Range("A$r").Value = $data[$s]["A"]
    $ws.Range("B$r").Value = $data[$s]["B"]
    $ws.Range("E$r").Value = $data[$s]["E"]
    $ws.Range("F$r").Value = $data[$s]["F"]
    $ws.Range("G$r").Value = $data[$s]["G"]
    $ws.Range("H$r").Value = $data[$s]["H"]
    $ws.Range("P$r").Value = $data[$s]["P"]
    $ws.Range("Q$r").Value = $data[$s]["Q"]
    $ws.Range("R$r").Value = $data[$s]["R"]
    $ws.Range("S$r").Value = $data[$s]["S"]
    $ws.Range("AW$r").Value = $data[$s]["AW"]
    $ws.Range("AX$r").Value = $data[$s]["AX"]

    # Z / AB (start/end time) - only set when the source row actually had a
    # value, otherwise clear the destination cell entirely.
    if ($null -eq $data[$s]["Z"] -or $data[$s]["Z"] -eq "") {
        $ws.Range("Z$r").ClearContents()
    } else {
        $ws.Range("Z$r").Value = $data[$s]["Z"]
    }
    if ($null -eq $data[$s]["AB"] -or $data[$s]["AB"] -eq "") {
        $ws.Range("AB$r").ClearContents()
    } else {
        $ws.Range("AB$r").Value = $data[$s]["AB"]
    }

    # AC (public comment) - only present on the bird-sighting rows.
    if ($null -eq $data[$s]["AC"] -or $data[$s]["AC"] -eq "") {
        $ws.Range("AC$r").ClearContents()
    } else {
        $ws.Range("AC$r").Value = $data[$s]["AC"]
    }

    # AF (determination method) - present (but empty) on some rows only.
    # A plain "" assignment removes the cell entirely, so when the source
    # row had an (empty) text cell we instead write a lone apostrophe
    # (Excel's "force text, empty" entry) which yields a real empty string
    # cell, then reset the style back to Normal (the apostrophe nudges the
    # cell onto a quote-prefixed style otherwise).
    if ($null -eq $data[$s]["AF"]) {
        $ws.Range("AF$r").ClearContents()
    } else {
        $ws.Range("AF$r").Value = "'"
        $ws.Range("AF$r").Style = "Normal"
    }
}
